$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.803.94"
$ws.Range("E2").Value = "  +6.35%  "

$ws.Range("D3").Value = "2.410.21"
$ws.Range("E3").Value = "  +2.52%  "

$ws.Range("E4").Value = "  +0.55%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "615.69"
$ws.Range("E5").Value = "  +11.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.31"
$ws.Range("E6").Value = "  +6.90%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.50%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.539"
$ws.Range("E8").Value = "  +2.43%  "

$ws.Range("D9").Value = "2.436.47"
$ws.Range("E9").Value = "  +4.06%  "

$ws.Range("E10").Value = "  +6.62%  "

$ws.Range("E11").Value = "  +1.36%  "

$ws.Range("E12").Value = "  +4.71%  "

$ws.Range("E13").Value = "  +5.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.94"
$ws.Range("E14").Value = "  +7.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000178"
$ws.Range("E15").Value = "  +10.17%  "

$ws.Range("D16").Value = "2.939.52"
$ws.Range("E16").Value = "  +5.54%  "

$ws.Range("D17").Value = "62.503.55"
$ws.Range("E17").Value = "  +5.22%  "

$ws.Range("D18").Value = "2.455.30"
$ws.Range("E18").Value = "  +5.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.95"
$ws.Range("E19").Value = "  -0.84%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.93"
$ws.Range("E20").Value = "  +6.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.29"
$ws.Range("E21").Value = "  +2.48%  "

$ws.Range("E22").Value = "  +4.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.03"
$ws.Range("E23").Value = "  +15.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.72"
$ws.Range("E25").Value = "  +3.41%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "615.41"
$ws.Range("E26").Value = "  +14.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.35"
$ws.Range("E27").Value = "  +4.76%  "

$ws.Range("D28").Value = "0.0₃0981"
$ws.Range("E28").Value = "  +9.97%  "

$ws.Range("D29").Value = "2.565.05"
$ws.Range("E29").Value = "  +3.73%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("E32").Value = "  +10.84%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.137"
$ws.Range("E33").Value = "  +7.84%  "

$ws.Range("E34").Value = "  +5.73%  "

$ws.Range("E35").Value = "  +6.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.996"
$ws.Range("E36").Value = "  -0.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.75"
$ws.Range("E37").Value = "  +6.65%  "

$ws.Range("E38").Value = "  +3.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "152.35"
$ws.Range("E39").Value = "  +1.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.37"
$ws.Range("E40").Value = "  +8.70%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.60"
$ws.Range("E41").Value = "  +3.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.77"
$ws.Range("E42").Value = "  +21.11%  "

$ws.Range("E43").Value = "  +9.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.36"
$ws.Range("E44").Value = "  +3.17%  "

$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").Value = "0.0₆0284"
$ws.Range("E46").Value = "  +1.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "144.02"
$ws.Range("E47").Value = "  +5.20%  "

$ws.Range("E48").Value = "  +3.56%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.17"
$ws.Range("E49").Value = "  +7.89%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.599"
$ws.Range("E50").Value = "  +3.99%  "

$ws.Range("E51").Value = "  +4.51%  "
